# Generate Report for Handback
#
# Updates the localization-status workbook to reflect that the handback
# (translation round-trip) has now completed for both the zh-cn and the
# de-de targets:
#   - the "Status" column (shared across the Overview + per-locale sheets)
#     flips from "Ready for handoff" to "Handed back: in sync with en-US"
#   - each per-locale sheet grows two new populated columns:
#       F "Latest Target File"   -> the (localized) target file name, linked
#       G "Latest Handback File" -> the handed-back xlf file name, linked
#   - the "Latest Handback DateTime" column (H) moves on from the
#     "0001-01-01 00:00:00" placeholder to a real timestamp per locale

$wb = $excel.ActiveWorkbook

$mdFile1 = "2246d8cb-028b-463f-8a0a-0d8d45762021.md"
$mdFile2 = "ffff7e577d15-3540-4630-af39-f5b803a4b64e.md"
$zhHandbackXlf = "2246d8cb-028b-463f-8a0a-0d8d45762021.1c8bef3171a484bb26a51d73e8ed7926ea9ada6e.zh-cn.xlf"
$deHandbackXlf = "2246d8cb-028b-463f-8a0a-0d8d45762021.1c8bef3171a484bb26a51d73e8ed7926ea9ada6e.de-de.xlf"

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1) Overview sheet: Status columns (B = zh-cn, C = de-de) for both rows
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusHandedBack
$wsOverview.Range("C2").Value = $statusHandedBack
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack

# ---------------------------------------------------------------------
# 2) zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status column now reflects the completed handback
$wsZh.Range("C2").Value = $statusHandedBack
$wsZh.Range("C3").Value = $statusHandedBack

# New "Latest Target File" (F) / "Latest Handback File" (G) links
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5ef689c6b95751645654ade0123a3009a5fc13a9/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/hb/$mdFile1", "", "", $mdFile1)
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5ef689c6b95751645654ade0123a3009a5fc13a9/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/hb/$zhHandbackXlf", "", "", $zhHandbackXlf)
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5ef689c6b95751645654ade0123a3009a5fc13a9/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/hb/$mdFile2", "", "", $mdFile2)
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5ef689c6b95751645654ade0123a3009a5fc13a9/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/hb/$zhHandbackXlf", "", "", $zhHandbackXlf)

# Latest Handback DateTime (H) - zh-cn finished first
$wsZh.Range("H2").Value = "2016-03-13 11:03:58"
$wsZh.Range("H3").Value = "2016-03-13 11:03:58"

# ---------------------------------------------------------------------
# 3) de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Status column now reflects the completed handback
$wsDe.Range("C2").Value = $statusHandedBack
$wsDe.Range("C3").Value = $statusHandedBack

# New "Latest Target File" (F) / "Latest Handback File" (G) links
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9450365fc5bac35f81aa7dd0bb939e580de73229/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/hb/$mdFile1", "", "", $mdFile1)
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9450365fc5bac35f81aa7dd0bb939e580de73229/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/hb/$deHandbackXlf", "", "", $deHandbackXlf)
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9450365fc5bac35f81aa7dd0bb939e580de73229/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/hb/$mdFile2", "", "", $mdFile2)
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9450365fc5bac35f81aa7dd0bb939e580de73229/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/hb/$deHandbackXlf", "", "", $deHandbackXlf)

# Latest Handback DateTime (H) - de-de finished a little later
$wsDe.Range("H2").Value = "2016-03-13 11:04:05"
$wsDe.Range("H3").Value = "2016-03-13 11:04:05"
